$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 15 ("Tension Spring") so that the existing
# row 15 data shifts down to row 16, and the new "MR115 Bearings" row
# becomes the new row 15.
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with the MR115 Bearings entry.
$ws.Range("A15").Value = "MR115 Bearings"
$ws.Range("B15").Value = 2
$ws.Range("D15").Value = "Likely to be replaced"
